$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Move Robot2 to location (2, 8) and remove the toolkit."
$ws.Range("B2").Value = 154.496142
$ws.Range("C2").Value = 9970
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.0306"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "a5d35f80-ac52-44af-bcf5-e576ebca0d90"

$ws.Range("A3").Value = "Move Robot26 to location (11, 4) and remove the liquid spill."
$ws.Range("B3").Value = 97.85663
$ws.Range("C3").Value = 9858
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.0297"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "75a7a699-fc3e-49dd-8100-df997e9c2b24"

$ws.Range("A4").Value = "Move Robot42 to location (9, 5) and remove the large debris."
$ws.Range("B4").Value = 80.989754
$ws.Range("C4").Value = 10560
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.03432"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "5811f724-f5b9-4a2a-a686-a97500191ab5"

$ws.Range("A5").Value = "Move Robot48 to location (5, 6) and remove the dust."
$ws.Range("B5").Value = 80.783254
$ws.Range("C5").Value = 10389
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.03201"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "c2dff074-5717-4bd6-b356-73bab2117dae"

$ws.Range("A6").Value = "Move Robot31 to location (9, 4) and remove the grass."
$ws.Range("B6").Value = 70.487014
$ws.Range("C6").Value = 9939
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.03039"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "efb39e2b-48e7-444a-85a7-ed5f08d91b1d"

$ws.Range("A7").Value = "Move Robot8 to location (8, 12) and remove the small debris."
$ws.Range("B7").Value = 74.856099
$ws.Range("C7").Value = 10435
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.03243"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "66e2d32f-5fcc-4f45-bd04-956263fb1c74"

$ws.Range("A8").Value = "Move Robot23 to location (11, 1) and remove the vehicle."
$ws.Range("B8").Value = 86.521565
$ws.Range("C8").Value = 10323
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.03075"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "622dc99c-9f86-42a7-83a0-f858a24d0ef3"

$ws.Range("A9").Value = "Move Robot23 to location (12, 10) and remove the construction materials."
$ws.Range("B9").Value = 83.645961
$ws.Range("C9").Value = 10585
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0342"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "f8c189b3-80a2-459f-ad93-e1da968400fb"

$ws.Range("A10").Value = "Move Robot14 to location (7, 11) and remove the tree branches."
$ws.Range("B10").Value = 76.400476
$ws.Range("C10").Value = 10441
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03237"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "b9fdcafc-45ec-421e-a6ae-2fbea5564198"

$ws.Range("A11").Value = "Move Robot15 to location (5, 3) and remove the screws."
$ws.Range("B11").Value = 66.392388
$ws.Range("C11").Value = 10015
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03141"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "97beca98-1d4f-451f-904d-57c6c44acb27"
